$d = $word.ActiveDocument

# 1) Rename the first bullet.
$d.Content.Find.Execute("Jouer une partie", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Démarrer la partie", 2)

# 2) Locate that paragraph (the one that now reads "Démarrer la partie") and
#    append five new list paragraphs after it, each inheriting its
#    Paragraphedeliste / Tahoma formatting.
$pDemarrer = $d.Paragraphs(3)
$pDemarrer.Range.InsertParagraphAfter()
$pDemarrer.Range.InsertParagraphAfter()
$pDemarrer.Range.InsertParagraphAfter()
$pDemarrer.Range.InsertParagraphAfter()
$pDemarrer.Range.InsertParagraphAfter()

# 3) Fill in the plain (single-run) paragraphs.
$d.Paragraphs(4).Range.Text = "Préparer son tour"

# 4) "Préparation " is authored as two runs ("P" + "réparation ") - type the
#    whole string then force a run split after the first character by
#    toggling a character property on it.
$d.Paragraphs(5).Range.Text = "Préparation "
$pPrep = $d.Paragraphs(5)
$rSplit1 = $d.Range($pPrep.Range.Start, $pPrep.Range.Start + 1)
$rSplit1.Bold = 1
$rSplit1.Bold = 0

# 5) "Envoie de flotte" is likewise authored as "E" + "nvoie de flotte".
$d.Paragraphs(6).Range.Text = "Envoie de flotte"
$pEnvoi = $d.Paragraphs(6)
$rSplit2 = $d.Range($pEnvoi.Range.Start, $pEnvoi.Range.Start + 1)
$rSplit2.Bold = 1
$rSplit2.Bold = 0

# 6) "Exécuter le tour" - this paragraph carries the relocated _GoBack
#    bookmark (it used to sit at the end of "Jouer une partie").
$d.Paragraphs(7).Range.Text = "Exécuter le tour"

# 7) Final new paragraph.
$d.Paragraphs(8).Range.Text = "Circonstance de fin de partie"

# 8) Re-seat the _GoBack bookmark around "Exécuter le tour" - adding a
#    bookmark under a name that already exists replaces the old one, so the
#    stale bookmark that used to close out "Jouer une partie" is removed as
#    a side effect of this call.
$pExec = $d.Paragraphs(7)
$bmRange = $d.Range($pExec.Range.Start, $pExec.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
